# #5: cash & deposit done
# Fill in the bank/deposit ("存款") sheet with the full metadata columns
# (property_category, category, date, legislator_name, legislator_id,
# source_file, index) that the other sheets (building/car/stock/...) already
# carry, and fix the header row (which previously repeated data values
# instead of real column labels).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)   # "存款" (deposit) sheet

# ---------------------------------------------------------------------
# Header row (row 1): real column-name labels
# ---------------------------------------------------------------------
$ws.Range("B1").Value = "bank"
$ws.Range("C1").Value = "deposit_type"
$ws.Range("D1").Value = "currency"
$ws.Range("E1").Value = "owner"
$ws.Range("F1").Value = "total"
$ws.Range("G1").Value = "property_category"
$ws.Range("H1").Value = "category"
$ws.Range("I1").Value = "date"
$ws.Range("J1").Value = "legislator_name"
$ws.Range("K1").Value = "legislator_id"
$ws.Range("L1").Value = "source_file"
$ws.Range("M1").Value = "index"

# New header cells (G1:M1) need the same bold/bordered look the existing
# header cells (B1:F1) already have. (Style object copies don't stick in
# this engine, so reproduce the look property-by-property instead.)
$headerCols = 7, 8, 9, 10, 11, 12, 13
foreach ($col in $headerCols) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

# ---------------------------------------------------------------------
# Data rows (2-4): move "total" into column F and append the metadata
# columns that mirror every other sheet in this workbook.
# ---------------------------------------------------------------------
$totals = @{ 2 = 813536; 3 = 673160.62; 4 = 453751.6 }
$indexes = @{ 2 = 44; 3 = 45; 4 = 46 }

for ($r = 2; $r -le 4; $r++) {
    $ws.Cells.Item($r, 6).Value = $totals[$r]          # F: total
    $ws.Cells.Item($r, 7).Value = "deposit"             # G: property_category
    $ws.Cells.Item($r, 8).Value = "normal"              # H: category
    $ws.Cells.Item($r, 9).Value = "'2013-11-12"         # I: date (keep as text)
    # Typing a leading apostrophe marks the cell "text-quoted"; strip that
    # back off so the cell looks like an ordinary data cell (no special
    # style), matching how the other plain text columns come out.
    $ws.Cells.Item($r, 9).Style = $ws.Cells.Item($r, 8).Style
    $ws.Cells.Item($r, 10).Value = "王育敏"              # J: legislator_name
    $ws.Cells.Item($r, 11).Value = 1728                 # K: legislator_id
    $ws.Cells.Item($r, 12).Value = "tmped871"           # L: source_file
    $ws.Cells.Item($r, 13).Value = $indexes[$r]         # M: index
}
